# Update countries & provincias Spain
# Refresh the "Pais" COVID table with the new pull: a handful of country
# rows get updated totals, three countries swap rank with their neighbour
# (Israel/Canada, El Salvador/Australia, Eslovaquia/Congo) because their
# updated total-case counts change the sort order, and the "last updated"
# timestamp moves from 09:51 to 11:08.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row layout: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

$rows = @(
    @{ Row = 25;  Pais = "Filipinas";   B = 238727; C = 1383; D = 184906; E = 49931; F = 0; G = 15; H = 3890 },
    @{ Row = 28;  Pais = "Israel";      B = 131970; C = 1326; D = 103849; E = 27099; F = 0; G = 3;  H = 1022 },
    @{ Row = 29;  Pais = "Canada";      B = 131895; C = 0;    D = 116357; E = 6393;  F = 0; G = 0;  H = 9145 },
    @{ Row = 30;  Pais = "Bolivia";     B = 120769; C = 528;  D = 71489;  E = 42272; F = 0; G = 40; H = 7008 },
    @{ Row = 32;  Pais = "Ecuador";     B = 109784; C = 0;    D = 91242;  E = 8018;  F = 0; G = 0;  H = 10524 },
    @{ Row = 49;  Pais = "Polonia";     B = 71126;  C = 302;  D = 55113;  E = 13889; F = 0; G = 4;  H = 2124 },
    @{ Row = 71;  Pais = "Austria";     B = 29561;  C = 290;  D = 25300;  E = 3515;  F = 0; G = 10; H = 746 },
    @{ Row = 73;  Pais = "El Salvador"; B = 26413;  C = 105;  D = 16137;  E = 9512;  F = 0; G = 5;  H = 764 },
    @{ Row = 74;  Pais = "Australia";   B = 26322;  C = 43;   D = 22603;  E = 2957;  F = 0; G = 9;  H = 762 },
    @{ Row = 90;  Pais = "Croacia";     B = 12081;  C = 117;  D = 9266;   E = 2614;  F = 0; G = 3;  H = 201 },
    @{ Row = 113; Pais = "Hong Kong";   B = 4890;   C = 11;   D = 4524;   E = 268;   F = 0; G = 4;  H = 98 },
    @{ Row = 117; Pais = "Eslovaquia";  B = 4636;   C = 22;   D = 2836;   E = 1763;  F = 0; G = 0;  H = 37 },
    @{ Row = 118; Pais = "Congo";       B = 4628;   C = 0;    D = 1742;   E = 2784;  F = 0; G = 0;  H = 102 },
    @{ Row = 133; Pais = "Lituania";    B = 3100;   C = 17;   D = 1955;   E = 1059;  F = 0; G = 0;  H = 86 },
    @{ Row = 137; Pais = "Estonia";     B = 2532;   C = 16;   D = 2176;   E = 292;   F = 0; G = 0;  H = 64 },
    @{ Row = 179; Pais = "Islas Feroe"; B = 413;    C = 0;    D = 409;    E = 4;     F = 0; G = 0;  H = 0 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 1).Value = $r.Pais
    $ws.Cells.Item($i, 2).Value = $r.B
    $ws.Cells.Item($i, 3).Value = $r.C
    $ws.Cells.Item($i, 4).Value = $r.D
    $ws.Cells.Item($i, 5).Value = $r.E
    $ws.Cells.Item($i, 6).Value = $r.F
    $ws.Cells.Item($i, 7).Value = $r.G
    $ws.Cells.Item($i, 8).Value = $r.H
}

# Update the "last refreshed" banner.
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 11:08"
